$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data keeps every "Price"/"Volume(1h)" column as literal text
# (e.g. "49.10", "0.0002000", "-3.47%") rather than as numbers, so each of
# these numeric-looking strings must be written with the cell pre-formatted
# as Text ("@") - otherwise Excel would helpfully "clean up" the value
# (dropping trailing zeros, turning "49.10" into 49.1, parsing "-3.47%" into
# the number -0.0347, etc.) which would not match the source data.
$textValues = @(
    @{ Cell = "D2";  Value = "309.26" }
    @{ Cell = "E2";  Value = "-3.47%" }
    @{ Cell = "D3";  Value = "49.10" }
    @{ Cell = "E3";  Value = "0.34%" }
    @{ Cell = "E4";  Value = "-1.75%" }
    @{ Cell = "D5";  Value = "0.07776" }
    @{ Cell = "E5";  Value = "-4.04%" }
    @{ Cell = "D6";  Value = "4.501" }
    @{ Cell = "E6";  Value = "-2.30%" }
    @{ Cell = "D7";  Value = "1.386" }
    @{ Cell = "E7";  Value = "15.49%" }
    @{ Cell = "D8";  Value = "1.557" }
    @{ Cell = "E8";  Value = "-6.35%" }
    @{ Cell = "D10"; Value = "0.2004" }
    @{ Cell = "E10"; Value = "3.03%" }
    @{ Cell = "D11"; Value = "0.04692" }
    @{ Cell = "E11"; Value = "4.03%" }
    @{ Cell = "D12"; Value = "0.09422" }
    @{ Cell = "E12"; Value = "-1.02%" }
    @{ Cell = "E13"; Value = "-0.17%" }
    @{ Cell = "D14"; Value = "0.001268" }
    @{ Cell = "E14"; Value = "-4.53%" }
    @{ Cell = "D15"; Value = "0.04174" }
    @{ Cell = "E15"; Value = "-2.82%" }
    @{ Cell = "D16"; Value = "0.005828" }
    @{ Cell = "E16"; Value = "-1.99%" }
    @{ Cell = "E17"; Value = "2,016.63%" }
    @{ Cell = "E18"; Value = "-0.73%" }
    @{ Cell = "D19"; Value = "2.237" }
    @{ Cell = "E19"; Value = "-8.17%" }
    @{ Cell = "D20"; Value = "0.3452" }
    @{ Cell = "E20"; Value = "1.76%" }
    @{ Cell = "D21"; Value = "7.932" }
    @{ Cell = "E21"; Value = "-3.14%" }
    @{ Cell = "E22"; Value = "-5.76%" }
    @{ Cell = "E23"; Value = "-0.86%" }
    @{ Cell = "D24"; Value = "0.001269" }
    @{ Cell = "E24"; Value = "-3.24%" }
    @{ Cell = "D25"; Value = "0.004039" }
    @{ Cell = "E25"; Value = "-4.88%" }
    @{ Cell = "D26"; Value = "0.0001350" }
    @{ Cell = "E26"; Value = "-0.14%" }
    @{ Cell = "D38"; Value = "0.02605" }
    @{ Cell = "E38"; Value = "-2.81%" }
    @{ Cell = "E39"; Value = "5.58%" }
    @{ Cell = "E40"; Value = "69.87%" }
    @{ Cell = "D41"; Value = "0.007938" }
    @{ Cell = "E41"; Value = "3.09%" }
    @{ Cell = "D42"; Value = "0.1423" }
    @{ Cell = "E42"; Value = "-1.06%" }
    @{ Cell = "D43"; Value = "0.008439" }
    @{ Cell = "E43"; Value = "9.46%" }
    @{ Cell = "D44"; Value = "0.008318" }
    @{ Cell = "E44"; Value = "2.67%" }
    @{ Cell = "D45"; Value = "0.3104" }
    @{ Cell = "E45"; Value = "-2.82%" }
    @{ Cell = "D46"; Value = "0.00007044" }
    @{ Cell = "E46"; Value = "0.41%" }
    @{ Cell = "E47"; Value = "-0.21%" }
    @{ Cell = "D48"; Value = "0.05399" }
    @{ Cell = "E48"; Value = "-11.25%" }
    @{ Cell = "D49"; Value = "0.002620" }
    @{ Cell = "E49"; Value = "-34.65%" }
    @{ Cell = "D50"; Value = "0.00002100" }
    @{ Cell = "E50"; Value = "-0.21%" }
    @{ Cell = "D51"; Value = "0.0002000" }
    @{ Cell = "E51"; Value = "-0.21%" }
)

foreach ($entry in $textValues) {
    $rng = $ws.Range($entry.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $entry.Value
    $rng.NumberFormat = "General"
}

# Rows 11 and 12 swapped coin identities (name + link); these are plain
# text already, so no special text-forcing is needed.
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
